$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subject = "Hebben we EcoPro-700 nog op voorraad?"
$content = "Testmail #6: $subject"

$body = @"
Beste afzender,
Hartelijk dank voor uw interesse in onze EcoPro-700. Op dit moment hebben we nog voldoende EcoPro-700 op voorraad. U kunt deze direct bestellen via onze website of neem contact met ons op als u meer informatie wenst.
Met vriendelijke groet,
[Naam]  
E-mailassistent  
[Bedrijfsnaam]
"@

$ws.Range("A22").Value = $content
$ws.Range("B22").Value = $body
$ws.Range("C22").Value = $subject
$ws.Range("D22").Value = "mailmind.test@zohomail.eu"
$ws.Range("E22").Value = "Productinformatie"
$ws.Range("F22").Value = "2025-07-29 21:39:42"
$ws.Range("G22").Value = "Ja"
$ws.Range("H22").Value = "Nee"
$ws.Range("I22").Value = "Ja"
$ws.Range("J22").Value = "Nee"

$ws.Rows.Item(22).EntireRow.AutoFit()
